$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: some of the new date strings (day <= 12) are ambiguous as
# DD-MM-YYYY vs MM-DD-YYYY, so a plain .Value assignment gets
# auto-converted by Excel into a date serial number instead of staying
# literal text (which is what the diff shows, since the source cells
# are plain strings with no date formatting). For those rows we
# temporarily force a text number format, assign the value, then
# restore the "Normal" style so no stray formatting is left behind.

# Row 3: date separator change + numeric updates
$ws.Range("A3").Value = "28-07-2022"
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Row 4: date separator change only (ambiguous date -> force text)
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "01-08-2022"
$ws.Range("A4").Style = "Normal"

# Row 5: date separator change + numeric updates (ambiguous date -> force text)
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "04-08-2022"
$ws.Range("A5").Style = "Normal"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

# Row 6: date separator change only (ambiguous date -> force text)
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "08-08-2022"
$ws.Range("A6").Style = "Normal"

# Row 7: date separator change only (ambiguous date -> force text)
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "11-08-2022"
$ws.Range("A7").Style = "Normal"

# Row 8: date separator change only
$ws.Range("A8").Value = "15-08-2022"

# Row 9: date separator change only
$ws.Range("A9").Value = "18-08-2022"

# Row 10: date separator change only
$ws.Range("A10").Value = "22-08-2022"

# Row 11: date separator change only
$ws.Range("A11").Value = "25-08-2022"

# Row 12: date separator change only
$ws.Range("A12").Value = "29-08-2022"

# Row 13: date separator change only (ambiguous date -> force text)
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "01-09-2022"
$ws.Range("A13").Style = "Normal"

# Row 14: date separator change only (ambiguous date -> force text)
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "05-09-2022"
$ws.Range("A14").Style = "Normal"

# Row 15: date separator change only (ambiguous date -> force text)
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "08-09-2022"
$ws.Range("A15").Style = "Normal"

# Row 16: date separator change only (ambiguous date -> force text)
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "12-09-2022"
$ws.Range("A16").Style = "Normal"

# Row 17: date separator change only
$ws.Range("A17").Value = "15-09-2022"

# Row 18: date separator change only
$ws.Range("A18").Value = "19-09-2022"

# Row 19: date separator change only
$ws.Range("A19").Value = "22-09-2022"

# Row 20: date separator change only
$ws.Range("A20").Value = "26-09-2022"

# Row 21: date separator change only
$ws.Range("A21").Value = "29-09-2022"
